# Fix mojibake "Â±" (UTF-8 bytes of U+00B1 re-decoded as Latin-1 then
# re-encoded as UTF-8) back to the correct "±" plus/minus sign across all
# worksheets that contain it.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Â±", "±")
}
